# Update CHNR yearly financials with latest reported figures
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CHNR")

$ws.Range("H8").Value = 6100
$ws.Range("I8").Value = 2200
$ws.Range("J8").Value = 6400
$ws.Range("F9").Value = 4700
$ws.Range("J9").Value = 3200
$ws.Range("G14").Value = 2600
$ws.Range("G17").Value = 5300
$ws.Range("H17").Value = 6600
$ws.Range("I17").Value = 4600
$ws.Range("J17").Value = 6000
$ws.Range("F18").Value = 2200
$ws.Range("G18").Value = -4100
$ws.Range("J18").Value = 400
$ws.Range("F20").Value = -2800
$ws.Range("I20").Value = 6900
$ws.Range("J20").Value = 4900
$ws.Range("H21").Value = 4300
$ws.Range("I21").Value = 7300
$ws.Range("J21").Value = 7800
$ws.Range("I22").Value = 6700
$ws.Range("J22").Value = 4800
$ws.Range("E23").Value = -700
$ws.Range("F23").Value = -600
$ws.Range("G23").Value = -3900
$ws.Range("I23").Value = -2300
$ws.Range("J23").Value = 500
$ws.Range("E26").Value = -700
$ws.Range("G26").Value = -4900
$ws.Range("E27").Value = -700
$ws.Range("G27").Value = -4900
$ws.Range("I27").Value = -3400
$ws.Range("E29").Value = -2800
$ws.Range("F29").Value = -5400
$ws.Range("H29").Value = -49600
$ws.Range("I29").Value = -10000
$ws.Range("J29").Value = -9600
$ws.Range("F32").Value = 2800
$ws.Range("I32").Value = -6900
$ws.Range("J32").Value = -4900
$ws.Range("D33").Value = -4500
$ws.Range("F33").Value = -6200
$ws.Range("G33").Value = -6500
$ws.Range("H33").Value = -49800
$ws.Range("I33").Value = -13300
$ws.Range("J33").Value = -9800
$ws.Range("D35").Value = -4500
$ws.Range("F35").Value = -6200
$ws.Range("G35").Value = -6500
$ws.Range("H35").Value = -49800
$ws.Range("I35").Value = -13300
$ws.Range("J35").Value = -9800
$ws.Range("D41").Value = 2800
$ws.Range("E41").Value = 2900
$ws.Range("F41").Value = 6700
$ws.Range("G41").Value = 7200
$ws.Range("H41").Value = 11400
$ws.Range("I41").Value = 31300
$ws.Range("J41").Value = 20300
$ws.Range("H42").Value = 600
$ws.Range("I42").Value = 3500
$ws.Range("J42").Value = 3000
$ws.Range("D43").Value = 1600
$ws.Range("I43").Value = 12000
$ws.Range("J43").Value = 7900
$ws.Range("E44").Value = 1600
$ws.Range("H45").Value = 431400
$ws.Range("J45").Value = 5100
$ws.Range("D46").Value = 4400
$ws.Range("E46").Value = 5400
$ws.Range("F46").Value = 8500
$ws.Range("G46").Value = 9300
$ws.Range("H46").Value = 444800
$ws.Range("I46").Value = 51500
$ws.Range("J46").Value = 38000
$ws.Range("D48").Value = 100
$ws.Range("E48").Value = 8100
$ws.Range("F48").Value = 7300
$ws.Range("H48").Value = 3200
$ws.Range("I48").Value = 348800
$ws.Range("J48").Value = 304700
$ws.Range("I52").Value = 23300
$ws.Range("J52").Value = 16200
$ws.Range("D54").Value = 4400
$ws.Range("E54").Value = 14100
$ws.Range("F54").Value = 16500
$ws.Range("G54").Value = 11900
$ws.Range("H54").Value = 448900
$ws.Range("I54").Value = 423600
$ws.Range("J54").Value = 358900
$ws.Range("E57").Value = 3000
$ws.Range("F57").Value = 3500
$ws.Range("I57").Value = 15500
$ws.Range("J57").Value = 16900
$ws.Range("I58").Value = 65900
$ws.Range("J58").Value = 33400
$ws.Range("D59").Value = 6200
$ws.Range("E59").Value = 8300
$ws.Range("F59").Value = 14900
$ws.Range("G59").Value = 3000
$ws.Range("H59").Value = 381000
$ws.Range("I59").Value = 81600
$ws.Range("J59").Value = 36500
$ws.Range("D60").Value = 6700
$ws.Range("E60").Value = 11300
$ws.Range("F60").Value = 18400
$ws.Range("G60").Value = 5600
$ws.Range("H60").Value = 381100
$ws.Range("I60").Value = 162900
$ws.Range("J60").Value = 86800
$ws.Range("I61").Value = 84500
$ws.Range("J61").Value = 77700
$ws.Range("G62").Value = 2900
$ws.Range("H62").Value = 29800
$ws.Range("I62").Value = 95000
$ws.Range("J62").Value = 101400
$ws.Range("D66").Value = 6700
$ws.Range("E66").Value = 12100
$ws.Range("F66").Value = 19100
$ws.Range("G66").Value = 8500
$ws.Range("H66").Value = 424800
$ws.Range("I66").Value = 356400
$ws.Range("J66").Value = 278700
$ws.Range("D72").Value = -48100
$ws.Range("E72").Value = -53000
$ws.Range("F72").Value = -57900
$ws.Range("G72").Value = -50900
$ws.Range("H72").Value = -30200
$ws.Range("I72").Value = 13900
$ws.Range("J72").Value = 31800
$ws.Range("E76").Value = 2000
$ws.Range("G76").Value = 3400
$ws.Range("H76").Value = 24000
$ws.Range("I76").Value = 67200
$ws.Range("J76").Value = 80200
$ws.Range("D81").Value = -4500
$ws.Range("F81").Value = -6200
$ws.Range("G81").Value = -6500
$ws.Range("H81").Value = -49800
$ws.Range("I81").Value = -13300
$ws.Range("J81").Value = -9800
$ws.Range("G83").Value = 1400
$ws.Range("H83").Value = 4100
$ws.Range("D89").Value = -2200
$ws.Range("E89").Value = -4200
$ws.Range("F89").Value = -3200
$ws.Range("G89").Value = -8900
$ws.Range("H89").Value = -22200
$ws.Range("I89").Value = -12000
$ws.Range("J89").Value = -9200
$ws.Range("F91").Value = -1600
$ws.Range("G91").Value = -8000
$ws.Range("H91").Value = -42800
$ws.Range("I91").Value = -52000
$ws.Range("J91").Value = -49700
$ws.Range("G94").Value = -5900
$ws.Range("H94").Value = -44000
$ws.Range("I94").Value = -54200
$ws.Range("J94").Value = -73400
$ws.Range("D100").Value = 1900
$ws.Range("G100").Value = -11200
$ws.Range("H100").Value = 68600
$ws.Range("I100").Value = 77300
$ws.Range("J100").Value = 84100
$ws.Range("E102").Value = -3900
$ws.Range("G102").Value = -26000
$ws.Range("H102").Value = 1900
$ws.Range("I102").Value = 11000
$ws.Range("J102").Value = 1500
